$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select the cell being edited, then update its value (Qualification -> Quali.)
$ws.Range("C1").Select()
$ws.Range("C1").Value = "Quali."
